$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.148.56'
$ws.Range('E2').Value = '  +0.71%  '
$ws.Range('D3').Value = '3.512.15'
$ws.Range('E3').Value = '  +0.19%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '595.72'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '173.60'
$ws.Range('E6').Value = '  +2.00%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.64%  '
$ws.Range('E9').Value = '  +5.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.29'
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('E11').Value = '  -1.00%  '
$ws.Range('D12').Value = '4.120.27'
$ws.Range('E12').Value = '  +0.23%  '
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.16'
$ws.Range('E14').Value = '  +2.65%  '
$ws.Range('D15').Value = '67.111.39'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '3.536.82'
$ws.Range('E17').Value = '  +0.47%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.35'
$ws.Range('E18').Value = '  -0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.19'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '396.52'
$ws.Range('E20').Value = '  +1.34%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.07'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.17'
$ws.Range('E22').Value = '  +0.03%  '
$ws.Range('E23').Value = '  +0.05%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.539'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.66'
$ws.Range('E25').Value = '  -3.74%  '
$ws.Range('E26').Value = '  -0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.23'
$ws.Range('E27').Value = '  -0.84%  '
$ws.Range('E28').Value = '  +0.27%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('E30').Value = '  -0.76%  '
$ws.Range('E31').Value = '  -2.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.07'
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '23.88'
$ws.Range('E33').Value = '  +1.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '7.39'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('E35').Value = '  +3.04%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '163.77'
$ws.Range('E36').Value = '  +0.65%  '
$ws.Range('E37').Value = '  +0.26%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.92'
$ws.Range('E38').Value = '  +0.51%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.08'
$ws.Range('E39').Value = '  +3.66%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0751'
$ws.Range('E41').Value = '  +0.12%  '
$ws.Range('B42').Value = 'InjectiveProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '27.48'
$ws.Range('E42').Value = '  +2.20%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.52'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('D44').Value = '2.821.02'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  +2.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '42.81'
$ws.Range('E46').Value = '  -0.54%  '
$ws.Range('E47').Value = '  -2.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '341.86'
$ws.Range('E48').Value = '  -3.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.98'
$ws.Range('E49').Value = '  +3.61%  '
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.51'
$ws.Range('E51').Value = '  -0.59%  '
